$d = $word.ActiveDocument

# 1) "nuestro banco de datos" -> "nuestra base de datos" in the intro paragraph
$d.Content.Find.Execute("nuestro banco de datos", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "nuestra base de datos", 2)

# 2) Update the cached TIME field result near the end of the document
$d.Content.Find.Execute("1 de noviembre de 2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "24 de febrero de 2025", 2)
